# Update "想去人数" (attendee interest count) figures in column F across
# three worksheets: 展览 (sheet1), 演出 (sheet2) and 全部类型 (sheet4).
# 本地生活 (sheet3) has no data rows and is left untouched.

$wb = $excel.ActiveWorkbook

# --- 展览 ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 332
$ws1.Range("F4").Value  = 267
$ws1.Range("F5").Value  = 10
$ws1.Range("F6").Value  = 3170
$ws1.Range("F7").Value  = 2083
$ws1.Range("F9").Value  = 148
$ws1.Range("F10").Value = 1180
$ws1.Range("F11").Value = 213
$ws1.Range("F12").Value = 1047

# --- 演出 -----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 31

# --- 全部类型 -------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 332
$ws4.Range("F4").Value  = 267
$ws4.Range("F5").Value  = 10
$ws4.Range("F6").Value  = 3170
$ws4.Range("F7").Value  = 2083
$ws4.Range("F9").Value  = 31
$ws4.Range("F10").Value = 148
$ws4.Range("F11").Value = 1180
$ws4.Range("F12").Value = 213
$ws4.Range("F13").Value = 1047
